$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 237 - this shifts the existing rows 237-248
# down to 238-249 (carrying their values/formatting with them), matching
# the diff where old row 237 data reappears (unchanged) as new row 238, and
# so on down the chain, with the former last row (248) becoming row 249.
$ws.Rows.Item(237).Insert()

# Populate the brand-new row 237 with its data.
$ws.Cells.Item(237, 1).Value = 9
$ws.Cells.Item(237, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(237, 3).Value = "Metropolitana"
$ws.Cells.Item(237, 4).Value = 44746
$ws.Cells.Item(237, 5).Value = 13
$ws.Cells.Item(237, 6).Value = 100112003
$ws.Cells.Item(237, 7).Value = "Ajo"
$ws.Cells.Item(237, 8).Value = "Chino"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 610
$ws.Cells.Item(237, 11).Value = 18000
$ws.Cells.Item(237, 12).Value = 20000
$ws.Cells.Item(237, 13).Value = 19000
$ws.Cells.Item(237, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(237, 15).Value = "China"
$ws.Cells.Item(237, 16).Value = 1900
$ws.Cells.Item(237, 17).Value = 10
$ws.Cells.Item(237, 18).Value = "Hortaliza"
